# Update cryptos list figures (Price / Volume(1h)) to refreshed values.
# Leading "'" forces Excel to keep these as text (matching the source data,
# which stores prices/percentages as text, e.g. "1.180", "0.1700", etc.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.024.75"
$ws.Range("D3").Value = "'1.826.55"
$ws.Range("E3").Value = "  -0.30%  "
$ws.Range("E4").Value = "  -0.25%  "
$ws.Range("D5").Value = "'311.86"
$ws.Range("E5").Value = "  -1.19%  "
$ws.Range("E6").Value = "  -0.27%  "
$ws.Range("D7").Value = "'0.4347"
$ws.Range("E7").Value = "  +0.97%  "
$ws.Range("D8").Value = "'0.3678"
$ws.Range("E8").Value = "  -0.87%  "
$ws.Range("D9").Value = "'0.07269"
$ws.Range("E9").Value = "  +0.03%  "
$ws.Range("D10").Value = "'0.8451"
$ws.Range("E10").Value = "  -2.76%  "
$ws.Range("D11").Value = "'20.67"
$ws.Range("E11").Value = "  -2.51%  "
$ws.Range("D12").Value = "'1.828.43"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").Value = "'6.661"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'0.07061"
$ws.Range("E14").Value = "  -0.36%  "
$ws.Range("D15").Value = "'5.297"
$ws.Range("E15").Value = "  -1.38%  "
$ws.Range("D16").Value = "'89.59"
$ws.Range("E16").Value = "  +1.77%  "
$ws.Range("D17").Value = "'1.003"
$ws.Range("E17").Value = "  -0.40%  "
$ws.Range("D18").Value = "'0.000008785"
$ws.Range("E18").Value = "  -1.60%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("E19").Value = "  -0.30%  "
$ws.Range("E20").Value = "  -2.41%  "
$ws.Range("D21").Value = "'27.111.39"
$ws.Range("E21").Value = "  -1.27%  "
$ws.Range("D22").Value = "'5.148"
$ws.Range("E22").Value = "  -0.57%  "
$ws.Range("D23").Value = "'10.89"
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").Value = "'2.056.67"
$ws.Range("E24").Value = "  -0.18%  "
$ws.Range("D25").Value = "'1.993"
$ws.Range("E25").Value = "  -1.22%  "
$ws.Range("E26").Value = "  -1.30%  "
$ws.Range("D27").Value = "'2.216"
$ws.Range("E27").Value = "  +2.85%  "
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("D29").Value = "'5.233"
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("D30").Value = "'116.98"
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("D31").Value = "'0.08757"
$ws.Range("E31").Value = "  -1.08%  "
$ws.Range("D32").Value = "'1.180"
$ws.Range("E32").Value = "  -2.55%  "
$ws.Range("D33").Value = "'0.7423"
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("D34").Value = "'2.906"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").Value = "'4.438"
$ws.Range("E35").Value = "  -1.61%  "
$ws.Range("D36").Value = "'1.000"
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("D37").Value = "'1.097"
$ws.Range("E37").Value = "  -2.40%  "
$ws.Range("D38").Value = "'0.01950"
$ws.Range("E38").Value = "  -0.89%  "
$ws.Range("D39").Value = "'0.05247"
$ws.Range("E39").Value = "  -0.91%  "
$ws.Range("D40").Value = "'7.239"
$ws.Range("E40").Value = "  +0.67%  "
$ws.Range("D41").Value = "'2.867"
$ws.Range("E41").Value = "  -0.85%  "
$ws.Range("D42").Value = "'0.1700"
$ws.Range("E42").Value = "  +0.96%  "
$ws.Range("D43").Value = "'0.5137"
$ws.Range("E43").Value = "  +0.89%  "
$ws.Range("D44").Value = "'8.580"
$ws.Range("E44").Value = "  -1.55%  "
$ws.Range("D45").Value = "'10.63"
$ws.Range("E45").Value = "  -0.41%  "
$ws.Range("D46").Value = "'0.4771"
$ws.Range("E46").Value = "  +0.55%  "
$ws.Range("D47").Value = "'106.04"
$ws.Range("E47").Value = "  -0.55%  "
$ws.Range("D48").Value = "'1.935"
$ws.Range("E48").Value = "  +5.42%  "
$ws.Range("E49").Value = "  -0.36%  "
$ws.Range("D50").Value = "'1.662"
$ws.Range("E50").Value = "  -0.90%  "
$ws.Range("D51").Value = "'0.06334"
